# Update the "Ponderacion_nueva" column (C2:C42) with recalculated weights.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 155.5511967810478
    3  = 11.44584116681023
    4  = 8.947694026983196
    5  = 14.07469242834858
    6  = 33.87009011600546
    7  = 10.57669567949825
    8  = 7.193157342316028
    9  = 26.09430253623052
    10 = 46.86961191603948
    11 = 8.879018979982334
    12 = 3.438921439602273
    13 = 6.757476936611639
    14 = 1.552942191858188
    15 = 1.348393933586805
    16 = 20.86023013685304
    17 = 21.6016252679161
    18 = 9.575369187743757
    19 = 1.070739980120957
    20 = 22.52172320945452
    21 = 65.91032252805249
    22 = 7.399920924684213
    23 = 2.37261210767492
    24 = 23.89891635629975
    25 = 5.31530094959355
    26 = 9.534754912635721
    27 = 19.74075302660244
    28 = 7.244109796542473
    29 = 5.196411889731843
    30 = 3.070439198167544
    31 = 1.740506298720755
    32 = 4.71125591253221
    33 = 2.213847214070778
    34 = 88.13740763263235
    35 = 8.098486456542437
    36 = 23.8317181920301
    37 = 4.079150103578045
    38 = 3.647161904701659
    39 = 8.554104779118044
    40 = 0.9666197475712638
    41 = 5.482188698219299
    42 = 286.624288113289
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
